$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the data range) used to stage each new value as a
# text formula result, then paste-special as values onto the real cell.
# This guarantees the target cell ends up holding a plain text value (not
# an auto-converted number) without leaving any unused cell style behind.
$helper = $ws.Range("Z1")

$helper.Formula = '="34.602.45"'
$helper.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$helper.Formula = '="  +0.24%  "'
$helper.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$helper.Formula = '="1.811.68"'
$helper.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$helper.Formula = '="  -0.04%  "'
$helper.Copy()
$ws.Range("E3").PasteSpecial(-4163)
$helper.Formula = '="  -0.14%  "'
$helper.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$helper.Formula = '="226.24"'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Formula = '="  -1.07%  "'
$helper.Copy()
$ws.Range("E5").PasteSpecial(-4163)
$helper.Formula = '="0.598"'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Formula = '="  +3.59%  "'
$helper.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$helper.Formula = '="  -0.11%  "'
$helper.Copy()
$ws.Range("E7").PasteSpecial(-4163)
$helper.Formula = '="37.88"'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Formula = '="  +8.26%  "'
$helper.Copy()
$ws.Range("E8").PasteSpecial(-4163)
$helper.Formula = '="  -3.00%  "'
$helper.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$helper.Formula = '="  -2.08%  "'
$helper.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$helper.Formula = '="0.0969"'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Formula = '="  +1.42%  "'
$helper.Copy()
$ws.Range("E11").PasteSpecial(-4163)
$helper.Formula = '="2.072.23"'
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$helper.Formula = '="  -0.07%  "'
$helper.Copy()
$ws.Range("E12").PasteSpecial(-4163)
$helper.Formula = '="11.40"'
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$helper.Formula = '="  +1.26%  "'
$helper.Copy()
$ws.Range("E13").PasteSpecial(-4163)
$helper.Formula = '="1.818.47"'
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$helper.Formula = '="  +0.34%  "'
$helper.Copy()
$ws.Range("E14").PasteSpecial(-4163)
$helper.Formula = '="0.635"'
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$helper.Formula = '="  -1.92%  "'
$helper.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$helper.Formula = '="34.538.81"'
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Formula = '="  +0.07%  "'
$helper.Copy()
$ws.Range("E16").PasteSpecial(-4163)
$helper.Formula = '="  -0.39%  "'
$helper.Copy()
$ws.Range("E17").PasteSpecial(-4163)
$helper.Formula = '="68.67"'
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$helper.Formula = '="  -0.78%  "'
$helper.Copy()
$ws.Range("E18").PasteSpecial(-4163)
$helper.Formula = '="244.65"'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$helper.Formula = '="  -0.36%  "'
$helper.Copy()
$ws.Range("E19").PasteSpecial(-4163)
$helper.Formula = '="0.0₃0777"'
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Formula = '="  -2.70%  "'
$helper.Copy()
$ws.Range("E20").PasteSpecial(-4163)
$helper.Formula = '="11.26"'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$helper.Formula = '="  -1.67%  "'
$helper.Copy()
$ws.Range("E21").PasteSpecial(-4163)
$helper.Formula = '="  -0.02%  "'
$helper.Copy()
$ws.Range("E22").PasteSpecial(-4163)
$helper.Formula = '="4.15"'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Formula = '="  -0.38%  "'
$helper.Copy()
$ws.Range("E23").PasteSpecial(-4163)
$helper.Formula = '="2.21"'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Formula = '="  +4.44%  "'
$helper.Copy()
$ws.Range("E24").PasteSpecial(-4163)
$helper.Formula = '="172.17"'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Formula = '="  -0.29%  "'
$helper.Copy()
$ws.Range("E25").PasteSpecial(-4163)
$helper.Formula = '="7.89"'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$helper.Formula = '="  -0.79%  "'
$helper.Copy()
$ws.Range("E26").PasteSpecial(-4163)
$helper.Formula = '="17.42"'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$helper.Formula = '="  +3.67%  "'
$helper.Copy()
$ws.Range("E27").PasteSpecial(-4163)
$helper.Formula = '="  -0.14%  "'
$helper.Copy()
$ws.Range("E29").PasteSpecial(-4163)
$helper.Formula = '="3.95"'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Formula = '="  -2.07%  "'
$helper.Copy()
$ws.Range("E30").PasteSpecial(-4163)
$helper.Formula = '="3.83"'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$helper.Formula = '="  -0.98%  "'
$helper.Copy()
$ws.Range("E31").PasteSpecial(-4163)
$helper.Formula = '="  -0.98%  "'
$helper.Copy()
$ws.Range("E32").PasteSpecial(-4163)
$helper.Formula = '="0.0522"'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$helper.Formula = '="  -2.16%  "'
$helper.Copy()
$ws.Range("E33").PasteSpecial(-4163)
$helper.Formula = '="  -0.56%  "'
$helper.Copy()
$ws.Range("E34").PasteSpecial(-4163)
$helper.Formula = '="1.368.08"'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Formula = '="  -1.94%  "'
$helper.Copy()
$ws.Range("E35").PasteSpecial(-4163)
$helper.Formula = '="0.657"'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$helper.Formula = '="  -3.53%  "'
$helper.Copy()
$ws.Range("E36").PasteSpecial(-4163)
$helper.Formula = '="  +1.36%  "'
$helper.Copy()
$ws.Range("E37").PasteSpecial(-4163)
$helper.Formula = '="  -4.08%  "'
$helper.Copy()
$ws.Range("E38").PasteSpecial(-4163)
$helper.Formula = '="  -1.29%  "'
$helper.Copy()
$ws.Range("E39").PasteSpecial(-4163)
$helper.Formula = '="  +8.34%  "'
$helper.Copy()
$ws.Range("E40").PasteSpecial(-4163)
$helper.Formula = '="  +1.38%  "'
$helper.Copy()
$ws.Range("E41").PasteSpecial(-4163)
$helper.Formula = '="Aave"'
$helper.Copy()
$ws.Range("B42").PasteSpecial(-4163)
$helper.Formula = '="https://coinranking.com/coin/ixgUfzmLR+aave-aave"'
$helper.Copy()
$ws.Range("C42").PasteSpecial(-4163)
$helper.Formula = '="81.23"'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Formula = '="  -3.09%  "'
$helper.Copy()
$ws.Range("E42").PasteSpecial(-4163)
$helper.Formula = '="ARBITRUM"'
$helper.Copy()
$ws.Range("B43").PasteSpecial(-4163)
$helper.Formula = '="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"'
$helper.Copy()
$ws.Range("C43").PasteSpecial(-4163)
$helper.Formula = '="0.941"'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Formula = '="  -2.29%  "'
$helper.Copy()
$ws.Range("E43").PasteSpecial(-4163)
$helper.Formula = '="MXToken"'
$helper.Copy()
$ws.Range("B44").PasteSpecial(-4163)
$helper.Formula = '="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"'
$helper.Copy()
$ws.Range("C44").PasteSpecial(-4163)
$helper.Formula = '="2.78"'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Formula = '="  -1.92%  "'
$helper.Copy()
$ws.Range("E44").PasteSpecial(-4163)
$helper.Formula = '="13.99"'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$helper.Formula = '="  +4.56%  "'
$helper.Copy()
$ws.Range("E45").PasteSpecial(-4163)
$helper.Formula = '="0.0501"'
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$helper.Formula = '="  -2.53%  "'
$helper.Copy()
$ws.Range("E46").PasteSpecial(-4163)
$helper.Formula = '="1.972.57"'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Formula = '="5.84"'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$helper.Formula = '="  -2.55%  "'
$helper.Copy()
$ws.Range("E48").PasteSpecial(-4163)
$helper.Formula = '="  -0.06%  "'
$helper.Copy()
$ws.Range("E49").PasteSpecial(-4163)
$helper.Formula = '="103.24"'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$helper.Formula = '="  -1.90%  "'
$helper.Copy()
$ws.Range("E50").PasteSpecial(-4163)
$helper.Formula = '="  -7.47%  "'
$helper.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$helper.ClearContents()
$excel.CutCopyMode = 0
